# Added v0 11v0 and 12v0
# T2 no PBS chase 1 and 10 mL/hr scripts

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("v0")

# Label the existing "v0-protocol-7v3.json" row (row 13) as the PBS-chase variant.
$ws.Range("B13").Value2 = "T2 (PBS Chase) "

# --- New row 17: v0-protocol-11v0.json / T2 (No PBS Chase) 1 mL/hr ---
$ws.Range("A17").Value2 = "v0-protocol-11v0.json"
$ws.Range("C17").Value2 = "5 mL"
$ws.Range("D17").Value2 = "5 mL"
$ws.Range("E17").Value2 = "1 hour"
$ws.Range("F17").Value2 = "1.0 mL"
$ws.Range("G17").Value2 = "1 mL/hr"
$ws.Range("G17").Font.Bold = $true
$ws.Range("G17").Font.Color = 255

# --- New row 18: v0-protocol-12v0.json / T2 (No PBS Chase) 10 mL/hr ---
$ws.Range("A18").Value2 = "v0-protocol-12v0.json"
$ws.Range("C18").Value2 = "5 mL"
$ws.Range("D18").Value2 = "5 mL"
$ws.Range("E18").Value2 = "1 hour"
$ws.Range("F18").Value2 = "1.0 mL"
$ws.Range("G18").Value2 = "10 mL/hr"
$ws.Range("G18").Font.Bold = $true
$ws.Range("G18").Font.Color = 255

# Columns H:L are not applicable for these scripts.
$ws.Range("H17:L18").ClearFormats()
$ws.Range("H17:L18").Value2 = "N/A"

# Notes column.
$ws.Range("M18").Value2 = "Exact copy of 11v0 except 10 mL/hr instead of 1 mL/hr"
$ws.Range("M17").Value2 = "Can use 7v3 as reference - remove extra 2 min incubation after F-127; change sample flow rate to 1 mL/hr; sample pull 500 ul to waste syringe then 1 mL to lysate/collection syringe"

# Name column.
$ws.Range("B17").Value2 = "T2 (No PBS Chase) 1 mL/hr"
$ws.Range("B18").Value2 = "T2 (No PBS Chase) 10 mL/hr"

# Reflect the user's last selection/scroll position on the sheet.
$ws.Activate() | Out-Null
$ws.Range("B19").Select() | Out-Null
